$wb = $excel.ActiveWorkbook

# 1. Rename the "Google" sheet to "TC" (sheetId 5 / r:id rId1, first tab).
$wb.Worksheets.Item("Google").Name = "TC"

# 2. On the "DragnDrop" sheet, push the current row 4 ("quit") down to
#    row 5, then turn row 4 into a new "refresh" step.
$ws = $wb.Worksheets.Item("DragnDrop")

# Duplicate row 4 (values + formats) into row 5, cell by cell, so existing
# cell styles are reused instead of cloning new style entries.
$cols = @("A","B","C","D","E","F","G","H","I","J")
foreach ($col in $cols) {
    $src = $ws.Range($col + "4")
    $dst = $ws.Range($col + "5")
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value2 = $src.Value2
}
$excel.CutCopyMode = 0

# The row that moved down is now step 4.
$ws.Range("A5").Value2 = 4

# Replace row 4 in place with the new step: only Step / Action / Screenshot
# columns stay populated.
$ws.Range("B4:H4").Clear()
$ws.Range("J4").Clear()

# Reuse the existing cell styles (instead of cloning new ones) for the
# cells that keep formatting in the new row 4.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A4").Value2 = 3
$ws.Range("B4").Value2 = "refresh"
$ws.Range("I4").Value2 = $true

# 3. Update the active selection on the sheet.
$ws.Activate()
$ws.Range("E9").Select()
